# Update a book successfully in NewBookWindow
# - Add a "Testing" task row (Member = "All")
# - Fill in the "Note" column (D) for the Coding sheet's CRUD tasks
# - Resize column D to fit the new multi-line notes
# - Make the "Coding" sheet the active/selected tab (was "Documenting")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Coding" sheet

# --- Row 8: new "Testing" task ---
$ws.Range("B8").Value = "Testing"
$ws.Range("C8").Value = "All"

# --- D5: Add a new book form note ---
$ws.Range("D5").Value = "- CRUD a book`n- input validation"
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = 28.5

# --- D4: Add a new member form note ---
$ws.Range("D4").Value = "- CRUD a member`n- input validation"
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = 28.5

# --- D7: Add a copy form note ---
$ws.Range("D7").Value = "- CRUD a book copy`n- input validation"
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Rows.Item(7).RowHeight = 28.5

# --- D3: Login form note ---
$ws.Range("D3").Value = "- login`n- show/hide corresponding menus"
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Rows.Item(3).RowHeight = 28.5

# --- widen column D to fit the new notes ---
$ws.Columns.Item(4).ColumnWidth = 32.59

# --- make Coding the active sheet/tab and set the new selection ---
$ws.Activate()
$ws.Range("D4").Select()
